$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (it holds values like "1.006" which
# Excel would otherwise auto-convert to numbers) without leaving a lingering
# number-format override: apply "@", write the values, then reset the style
# back to Normal so the cells end up stored exactly like the originals.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.717.11'
$ws.Range("D3").Value = '1.890.15'
$ws.Range("D4").Value = '1.006'
$ws.Range("D5").Value = '323.42'
$ws.Range("D6").Value = '1.004'
$ws.Range("D7").Value = '0.4526'
$ws.Range("D8").Value = '0.3793'
$ws.Range("D9").Value = '0.07714'
$ws.Range("D10").Value = '0.9697'
$ws.Range("D11").Value = '21.79'
$ws.Range("D12").Value = '1.906.15'
$ws.Range("D13").Value = '5.659'
$ws.Range("D14").Value = '6.910'
$ws.Range("D15").Value = '0.06974'
$ws.Range("D16").Value = '1.006'
$ws.Range("D17").Value = '84.28'
$ws.Range("D18").Value = '0.000009387'
$ws.Range("D19").Value = '16.55'
$ws.Range("D20").Value = '1.004'
$ws.Range("D21").Value = '28.731.83'
$ws.Range("D22").Value = '5.287'
$ws.Range("D23").Value = '11.15'
$ws.Range("D24").Value = '2.117.06'
$ws.Range("D25").Value = '2.058'
$ws.Range("D26").Value = '158.06'
$ws.Range("D27").Value = '18.88'
$ws.Range("D28").Value = '5.597'
$ws.Range("D29").Value = '117.46'
$ws.Range("D30").Value = '1.835'
$ws.Range("D31").Value = '0.09273'
$ws.Range("D32").Value = '0.8629'
$ws.Range("D34").Value = '1.230'
$ws.Range("D35").Value = '2.983'
$ws.Range("D36").Value = '0.05691'
$ws.Range("D38").Value = '1.004'
$ws.Range("D39").Value = '0.02027'
$ws.Range("D40").Value = '3.035'
$ws.Range("D41").Value = '7.484'
$ws.Range("D42").Value = '0.5477'
$ws.Range("D43").Value = '0.1749'
$ws.Range("D44").Value = '9.286'
$ws.Range("D45").Value = '0.000002738'
$ws.Range("D46").Value = '2.154'
$ws.Range("D47").Value = '0.5119'
$ws.Range("D48").Value = '0.06893'
$ws.Range("D49").Value = '11.05'
$ws.Range("D50").Value = '109.57'
$ws.Range("D51").Value = '1.748'

$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) / Coin / Link text updates
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("E40").Value = '  +10.47%  '
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  +4.98%  '
$ws.Range("E46").Value = '  +4.10%  '
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("E51").Value = '  -1.47%  '

Write-Output "Applied crypto list update"
